# Update "Berekening oversterfte" sheet with the latest CBS weekly figures.
# - revises several already-entered G (verwacht) / H (overleden) counts
# - turns column I ("oversterfte") into a real G-H formula (shared across I4:I25)
# - appends two new weeks (rows 24 and 25)
# - extends the totals row (28) to cover the new rows
# - moves the active selection from I24 to I6

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Revised weekly figures (G = verwacht aantal, H = overleden) ------------
$ws.Range("H3").Value  = 3052

$ws.Range("H4").Value  = 3093

$ws.Range("H5").Value  = 3114

$ws.Range("H6").Value  = 3114

$ws.Range("G7").Value  = 4977
$ws.Range("H7").Value  = 2908

$ws.Range("G8").Value  = 4299
$ws.Range("H8").Value  = 3009

$ws.Range("G9").Value  = 3905
$ws.Range("H9").Value  = 2925

$ws.Range("G10").Value = 3378
$ws.Range("H10").Value = 2978

$ws.Range("G11").Value = 2981
$ws.Range("H11").Value = 2932

$ws.Range("H12").Value = 3049

$ws.Range("G13").Value = 2768
$ws.Range("H13").Value = 2810

$ws.Range("G14").Value = 2724
$ws.Range("H14").Value = 2822

$ws.Range("G15").Value = 2681
$ws.Range("H15").Value = 2860

$ws.Range("G16").Value = 2691
$ws.Range("H16").Value = 2806

$ws.Range("G17").Value = 2690
$ws.Range("H17").Value = 2891

$ws.Range("G18").Value = 2659
$ws.Range("H18").Value = 3063

$ws.Range("G19").Value = 2633
$ws.Range("H19").Value = 2835

$ws.Range("G20").Value = 2610

$ws.Range("G21").Value = 2519
$ws.Range("H21").Value = 2848

$ws.Range("G22").Value = 2666
$ws.Range("H22").Value = 2835

$ws.Range("G23").Value = 2640

# --- Two newly-reported weeks -----------------------------------------------
$ws.Range("F24").Value = 32
$ws.Range("G24").Value = 2623
$ws.Range("H24").Value = 3103

$ws.Range("F25").Value = 33
$ws.Range("G25").Value = 3172
$ws.Range("H25").Value = 3114

# --- Column I: make it a real formula (G-H), shared down I4:I25 ------------
$ws.Range("I3").Formula = "=G3-H3"
$ws.Range("I4:I25").Formula = "=G4-H4"

# --- Extend the totals row to the new last data row (25) --------------------
$ws.Range("G28").Formula = "=SUM(G3:G25)"
$ws.Range("H28").Formula = "=SUM(H3:H25)"
$ws.Range("I28").Formula = "=SUM(I3:I25)"

# --- Move the selection / scroll position (was I24, now I6) ----------------
$ws.Range("I6").Select()
